# Apply "length_of_stay_categories" sheet restructuring:
#   - split the old "length_of_stay_simple" (1-13 / 14+) column into
#     "length_of_stay_simple_two_weeks" (same values) plus a brand-new
#     "length_of_stay_simple_week" (1-6 / 7+) column right after it
#   - keep "length_of_stay_three" (1-6 / 7-13 / 14+) as-is, just shifted
#   - split the old "length_of_stay_simple_order" column into
#     "length_of_stay_simple_two_weeks_order" (same values) plus a new
#     "length_of_stay_simple_week_order" (1 / 2) column right after it
#   - keep "length_of_stay_three_order" as-is, just shifted to the end

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("length_of_stay_categories")

# --- restructure columns -------------------------------------------------
# Before:  A length_of_stay | B length_of_stay_simple | C length_of_stay_three
#        | D length_of_stay_simple_order | E length_of_stay_three_order
# Insert a new blank column at C (week) and, after that, a new blank
# column at F (week_order). This shifts the old C/D/E columns to D/E/G
# (with F left blank for the new week_order column).
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(6).Insert()

# New layout:
#   A length_of_stay
#   B length_of_stay_simple_two_weeks      (was length_of_stay_simple)
#   C length_of_stay_simple_week           (new)
#   D length_of_stay_three                 (was C, unchanged)
#   E length_of_stay_simple_two_weeks_order(was length_of_stay_simple_order)
#   F length_of_stay_simple_week_order     (new)
#   G length_of_stay_three_order           (was E, unchanged)

# Give the new "week" column the same look (width / text format) as its
# neighbouring "simple" columns. The new "week_order" column (F) holds
# plain numbers like its neighbours E and G, so it needs no special
# formatting beyond what Insert() already gave it.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(3).NumberFormat = "@"

# --- header row (B/C first; E/F headers are set after the data below, to
# mirror the order the new shared strings were authored in) --------------
$ws.Cells.Item(1, 2).Value2 = "length_of_stay_simple_two_weeks"
$ws.Cells.Item(1, 3).Value2 = "length_of_stay_simple_week"

# --- data rows ---------------------------------------------------------
# row r (2..60) corresponds to length_of_stay day count A = r - 1 (1..59)
for ($r = 2; $r -le 60; $r++) {
    $day = $r - 1

    if ($day -le 6) {
        $weekSimple = "1-13"
        $weekOnly   = "1-6"
        $threeVal   = "1-6"
        $twoWeeksOrder = 1
        $weekOrder     = 1
        $threeOrder    = 1
    } elseif ($day -le 13) {
        $weekSimple = "1-13"
        $weekOnly   = "7+"
        $threeVal   = "7-13"
        $twoWeeksOrder = 1
        $weekOrder     = 2
        $threeOrder    = 2
    } else {
        $weekSimple = "14+"
        $weekOnly   = "7+"
        $threeVal   = "14+"
        $twoWeeksOrder = 2
        $weekOrder     = 2
        $threeOrder    = 3
    }

    $ws.Cells.Item($r, 2).Value2 = $weekSimple
    $ws.Cells.Item($r, 3).Value2 = $weekOnly
    $ws.Cells.Item($r, 4).Value2 = $threeVal
    $ws.Cells.Item($r, 5).Value2 = $twoWeeksOrder
    $ws.Cells.Item($r, 6).Value2 = $weekOrder
    $ws.Cells.Item($r, 7).Value2 = $threeOrder
}

$ws.Cells.Item(1, 5).Value2 = "length_of_stay_simple_two_weeks_order"
$ws.Cells.Item(1, 6).Value2 = "length_of_stay_simple_week_order"

# --- sheet view tidy-up (matches the saved state in the target file) ----
$ws.Activate()
$ws.Range("B1").Select() | Out-Null
